# mirroring.pptx -- "initial commit for site reorg"
#
#   1. The fixed Date placeholder ("4/30/21") is updated to "6/15/21" on
#      every slide layout, the slide master and the notes master.
#   2. The small version-tag callouts on the main diagram slide are
#      shortened from "v1.0"/"v2.0" to "v1"/"v2" (7 shapes total).

$p = $ppt.ActivePresentation

$oldDate = "4/30/21"
$newDate = "6/15/21"

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# 1a. Every slide layout's Date placeholder.
$master = $p.SlideMaster
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# 1b. The slide master's own Date placeholder.
Update-DatePlaceholder $master.Shapes

# 1c. The notes master's Date placeholder.
Update-DatePlaceholder $p.NotesMaster.Shapes

# 2. Version-tag shapes on the (single) diagram slide.
$slide = $p.Slides.Item(1)
$versionMap = @{
    "Rounded Rectangle 7"  = "v1"
    "Rounded Rectangle 17" = "v2"
    "Rounded Rectangle 52" = "v1"
    "Rounded Rectangle 53" = "v2"
    "Rounded Rectangle 75" = "v1"
    "Rounded Rectangle 76" = "v2"
    "Rounded Rectangle 63" = "v2"
}

foreach ($name in $versionMap.Keys) {
    $shp = $slide.Shapes.Item($name)
    $shp.TextFrame.TextRange.Text = $versionMap[$name]
}
